# apex reworked for 4 TCDS links
# Original sheet had one "rx" link row (row 5, index 0) and one "tx" link
# row (row 6, index 0), both against TCDS interface index 4/5. The rework
# expands this into 4 TCDS links (index 0..3) for each of rx and tx, and
# renumbers the underlying interface index to 28..31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing "rx" row (row 5) is renumbered too: interface index moves
# from 4 to 28 (the new base index for this link group).
$ws.Range("C5").Value = "28"

# Row 6 currently holds the single "tx" entry. Insert 3 blank rows above
# it (pushing it down to row 9) so we can add 3 more "rx" rows (6,7,8)
# that mirror the formatting of row 5 (the existing "rx" row above them).
$ws.Rows("6:8").Insert()

# Fill in the new rx rows (index 1,2,3), same direction/protocol as row 5,
# with interface index incrementing 29,30,31.
$ws.Range("A6").Value = "rx"
$ws.Range("B6").Value = "0"
$ws.Range("C6").Value = "29"
$ws.Range("D6").Value = "TCDS"
$ws.Range("E6").Value = "Q"
$ws.Range("F6").Value = "160.3144"
$ws.Range("G6").Value = "tcds_rx"
$ws.Range("H6").Value = "1"

$ws.Range("A7").Value = "rx"
$ws.Range("B7").Value = "0"
$ws.Range("C7").Value = "30"
$ws.Range("D7").Value = "TCDS"
$ws.Range("E7").Value = "Q"
$ws.Range("F7").Value = "160.3144"
$ws.Range("G7").Value = "tcds_rx"
$ws.Range("H7").Value = "2"

$ws.Range("A8").Value = "rx"
$ws.Range("B8").Value = "0"
$ws.Range("C8").Value = "31"
$ws.Range("D8").Value = "TCDS"
$ws.Range("E8").Value = "Q"
$ws.Range("F8").Value = "160.3144"
$ws.Range("G8").Value = "tcds_rx"
$ws.Range("H8").Value = "3"

# The original "tx" row (now shifted down to row 9) keeps direction tx /
# index 0, but its interface index moves from 5 to 28 to match the new
# numbering scheme.
$ws.Range("C9").Value = "28"

# Insert 3 more rows below row 9 for the remaining tx links (index 1,2,3),
# mirroring row 9's formatting.
$ws.Rows("10:12").Insert()

$ws.Range("A10").Value = "tx"
$ws.Range("B10").Value = "0"
$ws.Range("C10").Value = "29"
$ws.Range("D10").Value = "TCDS"
$ws.Range("E10").Value = "Q"
$ws.Range("F10").Value = "160.3144"
$ws.Range("G10").Value = "tcds_tx"
$ws.Range("H10").Value = "1"

$ws.Range("A11").Value = "tx"
$ws.Range("B11").Value = "0"
$ws.Range("C11").Value = "30"
$ws.Range("D11").Value = "TCDS"
$ws.Range("E11").Value = "Q"
$ws.Range("F11").Value = "160.3144"
$ws.Range("G11").Value = "tcds_tx"
$ws.Range("H11").Value = "2"

$ws.Range("A12").Value = "tx"
$ws.Range("B12").Value = "0"
$ws.Range("C12").Value = "31"
$ws.Range("D12").Value = "TCDS"
$ws.Range("E12").Value = "Q"
$ws.Range("F12").Value = "160.3144"
$ws.Range("G12").Value = "tcds_tx"
$ws.Range("H12").Value = "3"

# Move the active selection to H11, matching the edited workbook.
$ws.Range("H11").Select()
